$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: URL pattern changed from "/profile" to "/member"
$ws.Range("C6").Value = "/member"

# New row 7: additional signup endpoint entry
$ws.Range("C7").Value = "/member"
$ws.Range("D7").Value = "/signup"
$ws.Range("I7").Value = "회원가입"

# Remove the stray NoSuchElementException note that used to live at H11
$ws.Range("H11").ClearContents() | Out-Null

# Column H needs its own (wider) width now that it holds longer response names
$ws.Columns("H").ColumnWidth = 21.57

# Update the active selection to H4
$ws.Range("H4").Select() | Out-Null
